# Update crypto price/volume snapshot (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose ranking position is unchanged: only Price (D) / Volume(1h) (E) move ---
$ws.Range("D2").Value = "68.063.53"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "3.319.11"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'582.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").Value = "'175.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.25%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.583"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.11%  "
$ws.Range("D9").Value = "3.316.39"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("E10").Value = "  -3.70%  "
$ws.Range("D11").Value = "'0.578"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("E12").Value = "  -4.14%  "
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").Value = "'664.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.60%  "
$ws.Range("D15").Value = "3.861.45"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "68.211.84"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").Value = "3.325.88"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "'17.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("E21").Value = "  -2.22%  "
$ws.Range("D22").Value = "'0.890"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.25%  "

# --- Rows 23/24 swap ranking: Toncoin <-> InternetComputer(DFINITY) ---
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'17.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.02%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'5.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.44%  "

$ws.Range("D25").Value = "'97.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("E26").Value = "  -4.01%  "
$ws.Range("E27").Value = "  -5.85%  "
$ws.Range("D28").Value = "'9.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.30%  "
$ws.Range("D29").Value = "'33.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.88%  "
$ws.Range("D30").Value = "'8.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("D31").Value = "'7.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.96%  "
$ws.Range("D32").Value = "'590.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("D33").Value = "'10.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("D35").Value = "3.748.69"
$ws.Range("E35").Value = "  -6.31%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "'3.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.19%  "
$ws.Range("D38").Value = "'55.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  +0.86%  "

# --- Rows 40/41 swap ranking: Fetch.AI <-> InjectiveProtocol ---
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'32.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.66%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'2.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.11%  "

$ws.Range("D42").Value = "'3.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.29%  "
$ws.Range("D43").Value = "0.0₃0667"
$ws.Range("E43").Value = "  -5.20%  "
$ws.Range("D44").Value = "'0.333"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.74%  "
$ws.Range("D45").Value = "'3.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.00%  "
$ws.Range("D46").Value = "'0.0408"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("D47").Value = "'2.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("D48").Value = "'0.128"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").Value = "'1.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("D51").Value = "'129.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "
